$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# worldclim row: type column changes from imageCollection to image
$ws.Range("D3").Value = "image"

# corine row: geeSnippet column gets the /2018 suffix
$ws.Range("B5").Value = "COPERNICUS/CORINE/V20/100m/2018"

# restore the selected cell to B5 (was B14)
$ws.Range("B5").Select()
